# Generate Report for Handback
# Adds a new handed-back file (de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"  (File Name | Path And Name | Extension | Publish URL |
#                    zh-cn | de-de | Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Cells.Item(4, 1).Value = "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md"
$overview.Cells.Item(4, 3).Value = ".md"
$overview.Cells.Item(4, 5).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(4, 6).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(4, 7).Value = "2016-08-22 02:54:46"
$overview.Cells.Item(4, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$overview.Hyperlinks.Add(
    $overview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b16e0e6db7b536c1a9a5b6e40d2b2a8cca1f4b11/e2e/de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md",
    "",
    "",
    "e2e\de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md"
)
$overview.Range("B4").Font.Underline = $true
$overview.Range("B4").Font.Color = 6495237

$overviewTable = $overview.ListObjects.Item(1)
$overviewTable.Resize($overview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheets "zh-cn" and "de-de" (identical column layout)
#   Source File Name | File Extension | Status | Source Path | Priority |
#   Content Duplicate | Correspond Handoff File | Correspond Handoff Datetime |
#   Target File | Correspond Handback File | Correspond Handback DateTime |
#   Reference Tokens | To be localized | Dependency From | Has metadata |
#   Error Detail
# ---------------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; HandoffXlf = "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.63d3b4067c4d91fb4701731b8107fe798a18b2a4.zh-cn.xlf"; HandoffDate = "2016-08-22 02:54:39"; HandbackDate = "2016-08-22 02:55:13"; RepoSlug = "ol-test0-zhcn"; RepoCommit = "7d3f0f0f7c1d5b3b2f7c4f0a9c5b3f2e1d4c7b6a" },
    @{ Name = "de-de"; HandoffXlf = "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.63d3b4067c4d91fb4701731b8107fe798a18b2a4.de-de.xlf"; HandoffDate = "2016-08-22 02:54:46"; HandbackDate = "2016-08-22 02:55:20"; RepoSlug = "ol-test0-dede"; RepoCommit = "9a2c8e6f4b1d7c3a5e9f0b2d6c8a4e1f3b7d5c90" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    $ws.Cells.Item(4, 1).Value = "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md"
    $ws.Cells.Item(4, 2).Value = ".md"
    $ws.Cells.Item(4, 3).Value = "Handed back: in sync with en-US"
    $ws.Cells.Item(4, 4).Value = "e2e"
    $ws.Cells.Item(4, 5).Value = "ht"
    $ws.Cells.Item(4, 6).Value = "True"
    $ws.Cells.Item(4, 7).Value = $lang.HandoffXlf
    $ws.Cells.Item(4, 8).Value = $lang.HandoffDate
    $ws.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Cells.Item(4, 9).Value = "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md"
    $ws.Cells.Item(4, 10).Value = $lang.HandoffXlf
    $ws.Cells.Item(4, 11).Value = $lang.HandbackDate
    $ws.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Cells.Item(4, 13).Value = "True"
    $ws.Cells.Item(4, 15).Value = "False"

    $ws.Hyperlinks.Add(
        $ws.Range("A4"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b16e0e6db7b536c1a9a5b6e40d2b2a8cca1f4b11/e2e/de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md",
        "",
        "",
        "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md"
    )
    $ws.Range("A4").Font.Underline = $true
    $ws.Range("A4").Font.Color = 6495237

    $ws.Hyperlinks.Add(
        $ws.Range("I4"),
        "https://github.com/OpenLocalizationTestOrg/$($lang.RepoSlug)/blob/$($lang.RepoCommit)/e2e/de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md",
        "",
        "",
        "de3cccc0-1a06-4a66-9e0a-83c19c8ecd67.md"
    )
    $ws.Range("I4").Font.Underline = $true
    $ws.Range("I4").Font.Color = 6495237

    $lo = $ws.ListObjects.Item(1)
    $lo.Resize($ws.Range("A1:P4"))
}
